# Generate Report for Handback
# Row 7 ("e2dc9730-9dc9-4347-aac0-7bc4406ecb67") on both the "zh-cn" and
# "de-de" sheets moves from "Ready for handoff" into a handed-back-but-stale
# state: the Latest Target File / Latest Handback File / Latest Handback
# DateTime / Error Detail columns (I, J, K, P) get filled in, plus a new
# hyperlink on the Latest Target File cell.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4112b0f217fe0925d4ae81782eb51e98637e58fe/e2e/e2dc9730-9dc9-4347-aac0-7bc4406ecb67.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d369908de5fe5ff8652578e6070fadd9a8c28090/e2e/e2dc9730-9dc9-4347-aac0-7bc4406ecb67.md."

# ---- zh-cn sheet, row 7 ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$targetCellZh = $wsZh.Cells.Item(7, 9)
$wsZh.Hyperlinks.Add($targetCellZh, "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d369908de5fe5ff8652578e6070fadd9a8c28090/e2e/e2dc9730-9dc9-4347-aac0-7bc4406ecb67.md", "", "", "e2dc9730-9dc9-4347-aac0-7bc4406ecb67.md")

$wsZh.Cells.Item(7, 10).Value = "e2dc9730-9dc9-4347-aac0-7bc4406ecb67.1c72cced4175a16d6b9170eafc39c5af3a4c239b.zh-cn.xlf"
$wsZh.Cells.Item(7, 11).Value = "2016-08-23 16:59:32"
$wsZh.Cells.Item(7, 16).Value = $errorMessage

# ---- de-de sheet, row 7 ----
$wsDe = $wb.Worksheets.Item("de-de")

$targetCellDe = $wsDe.Cells.Item(7, 9)
$wsDe.Hyperlinks.Add($targetCellDe, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d369908de5fe5ff8652578e6070fadd9a8c28090/e2e/e2dc9730-9dc9-4347-aac0-7bc4406ecb67.md", "", "", "e2dc9730-9dc9-4347-aac0-7bc4406ecb67.md")

$wsDe.Cells.Item(7, 10).Value = "e2dc9730-9dc9-4347-aac0-7bc4406ecb67.1c72cced4175a16d6b9170eafc39c5af3a4c239b.de-de.xlf"
$wsDe.Cells.Item(7, 11).Value = "2016-08-23 16:59:39"
$wsDe.Cells.Item(7, 16).Value = $errorMessage
